# Remove the "عضاضه الجو" product row (row 25) entirely - Excel shifts all
# following rows (the "مبرد قدم" row, the totals row, and the footer row) up
# by one. Then refresh the computed total and the generation timestamp that
# live in the (now shifted) totals/footer rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Delete the whole row for "عضاضه الجو" (row 25).
$ws.Rows("25").Delete()

# 2) The totals row (old row 27) is now row 26; update its value to reflect
#    the removed 15.00 line item (935.12 -> 920.12).
$ws.Range("P26").Value = 920.12

# 3) The footer row (old row 28) is now row 27; refresh the "generated at"
#    timestamp shown there.
$ws.Range("A27").Value = "Tuesday, 9 September, 2025 12:49 PM"
